$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update quantities needed for the two charge-pump related parts.
$ws.Range("C5").Value = 7
$ws.Range("C18").Value = 17

# Update quantity ordered values and restyle those two cells.
$ws.Range("F5").Value = 10
$ws.Range("F5").Style = "Normal"
$ws.Range("F5").Font.Bold = $false

$ws.Range("F18").Value = 20
$ws.Range("F18").Style = "Normal"
$ws.Range("F18").Font.Bold = $false

# Update the view: clear the scrolled top-left cell and move the active selection.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("G13").Select()
